# Add team record columns (Wins / Losses / Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell onto the three new header cells, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record values for every data row (2-41).
$ws.Range("AD2:AD41").Value = 90
$ws.Range("AE2:AE41").Value = 72
$ws.Range("AF2:AF41").Value = 1
